$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.566.12'
$ws.Range("E2").Value = '  +2.58%  '

$ws.Range("D3").Value = '1.477.23'
$ws.Range("E3").Value = '  +3.89%  '

$ws.Range("D4").Value = '''1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.50%  '

$ws.Range("D5").Value = '''0.9591'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.78%  '

$ws.Range("D6").Value = '''276.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("D7").Value = '''0.3649'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.36%  '

$ws.Range("D8").Value = '''0.3054'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.65%  '

$ws.Range("D9").Value = '''39.66'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.05%  '

$ws.Range("D10").Value = '''1.056'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.07%  '

$ws.Range("D11").Value = '''0.06625'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.34%  '

$ws.Range("D12").Value = '''1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("D13").Value = '''18.24'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.31%  '

$ws.Range("D14").Value = '''5.473'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.75%  '

$ws.Range("D15").Value = '''6.172'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.50%  '

$ws.Range("D16").Value = '''0.00001030'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.75%  '

$ws.Range("D17").Value = '1.474.74'
$ws.Range("E17").Value = '  +3.63%  '

$ws.Range("D18").Value = '''0.05900'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.60%  '

$ws.Range("D19").Value = '''0.9646'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.23%  '

$ws.Range("D20").Value = '''69.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.21%  '

$ws.Range("D21").Value = '''5.471'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.56%  '

$ws.Range("D22").Value = '''14.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.51%  '

$ws.Range("D23").Value = '''11.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.65%  '

$ws.Range("D24").Value = '''2.247'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.77%  '

$ws.Range("D25").Value = '20.603.18'
$ws.Range("E25").Value = '  +2.50%  '

$ws.Range("D26").Value = '''142.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.95%  '

$ws.Range("D27").Value = '''2.131'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.88%  '

$ws.Range("D28").Value = '''17.23'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.72%  '

$ws.Range("D29").Value = '1.631.81'
$ws.Range("E29").Value = '  +3.21%  '

$ws.Range("D30").Value = '''113.63'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.40%  '

$ws.Range("D31").Value = '''3.914'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.71%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '''0.8196'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.73%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''4.984'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.47%  '

$ws.Range("D34").Value = '''0.07935'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.79%  '

$ws.Range("D35").Value = '''1.527'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.92%  '

$ws.Range("D36").Value = '''1.255'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +13.36%  '

$ws.Range("D37").Value = '''0.05757'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.60%  '

$ws.Range("D38").Value = '''4.736'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.03%  '

$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '''10.48'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.20%  '

$ws.Range("D40").Value = '''0.02041'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.11%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '''7.632'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.06%  '

$ws.Range("B42").Value = 'Frax'
$ws.Range("C42").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D42").Value = '''0.9583'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.80%  '

$ws.Range("D43").Value = '''0.1881'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.14%  '

$ws.Range("D44").Value = '''0.5293'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.81%  '

$ws.Range("D45").Value = '''3.504'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.17%  '

$ws.Range("D46").Value = '''12.07'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.75%  '

$ws.Range("D47").Value = '''117.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.61%  '

$ws.Range("D48").Value = '''0.5191'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.67%  '

$ws.Range("D49").Value = '''1.777'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.25%  '

$ws.Range("D50").Value = '''0.06466'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.66%  '

$ws.Range("D51").Value = '''0.9934'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.23%  '
